$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.07409833333333334
$ws.Cells.Item(2, 8).Value = 0.222295
$ws.Cells.Item(2, 9).Value = 0.02892330792240624
$ws.Cells.Item(2, 10).Value = 0.02892330792240624
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.05968133333333333
$ws.Cells.Item(2, 14).Value = 0.179044
$ws.Cells.Item(2, 15).Value = 0.02602747651633847
$ws.Cells.Item(2, 16).Value = 0.02602747651633848
$ws.Cells.Item(2, 17).Value = 0.004422287331111111
$ws.Cells.Item(2, 18).Value = 0.03980058598
$ws.Cells.Item(2, 19).Value = 0.000752800717725255
$ws.Cells.Item(2, 20).Value = 0.0007528007177252551

$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.07409833333333334
$ws.Cells.Item(3, 8).Value = 0.222295
$ws.Cells.Item(3, 9).Value = 0.02892330792240624
$ws.Cells.Item(3, 10).Value = 0.02892330792240624
$ws.Cells.Item(3, 15).Value = 0.144012433133819
$ws.Cells.Item(3, 16).Value = 0.144012433133819
$ws.Cells.Item(3, 17).Value = 0.02446892452944445
$ws.Cells.Item(3, 18).Value = 0.220220320765
$ws.Cells.Item(3, 19).Value = 0.004165315948184387
$ws.Cells.Item(3, 20).Value = 0.004165315948184386

$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.07409833333333334
$ws.Cells.Item(4, 8).Value = 0.222295
$ws.Cells.Item(4, 9).Value = 0.02892330792240624
$ws.Cells.Item(4, 10).Value = 0.02892330792240624
$ws.Cells.Item(4, 15).Value = 0.8299600903498424
$ws.Cells.Item(4, 16).Value = 0.8299600903498425
$ws.Cells.Item(4, 17).Value = 0.1410172050516667
$ws.Cells.Item(4, 18).Value = 1.269154845465
$ws.Cells.Item(4, 19).Value = 0.0240051912564966
$ws.Cells.Item(4, 20).Value = 0.0240051912564966

$ws.Cells.Item(5, 9).Value = 0.1644833827109413
$ws.Cells.Item(5, 10).Value = 0.1644833827109413
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.05968133333333333
$ws.Cells.Item(5, 14).Value = 0.179044
$ws.Cells.Item(5, 15).Value = 0.02602747651633847
$ws.Cells.Item(5, 16).Value = 0.02602747651633848
$ws.Cells.Item(5, 17).Value = 0.02514901758444444
$ws.Cells.Item(5, 18).Value = 0.22634115826
$ws.Cells.Item(5, 19).Value = 0.004281087380836937
$ws.Cells.Item(5, 20).Value = 0.004281087380836938

$ws.Cells.Item(6, 9).Value = 0.1644833827109413
$ws.Cells.Item(6, 10).Value = 0.1644833827109413
$ws.Cells.Item(6, 15).Value = 0.144012433133819
$ws.Cells.Item(6, 16).Value = 0.144012433133819
$ws.Cells.Item(6, 19).Value = 0.02368765215428379
$ws.Cells.Item(6, 20).Value = 0.02368765215428379

$ws.Cells.Item(7, 9).Value = 0.1644833827109413
$ws.Cells.Item(7, 10).Value = 0.1644833827109413
$ws.Cells.Item(7, 15).Value = 0.8299600903498424
$ws.Cells.Item(7, 16).Value = 0.8299600903498425
$ws.Cells.Item(7, 19).Value = 0.1365146431758205
$ws.Cells.Item(7, 20).Value = 0.1365146431758205

$ws.Cells.Item(8, 9).Value = 0.8065933093666526
$ws.Cells.Item(8, 10).Value = 0.8065933093666525
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.05968133333333333
$ws.Cells.Item(8, 14).Value = 0.179044
$ws.Cells.Item(8, 15).Value = 0.02602747651633847
$ws.Cells.Item(8, 16).Value = 0.02602747651633848
$ws.Cells.Item(8, 17).Value = 0.1233257061377778
$ws.Cells.Item(8, 18).Value = 1.10993135524
$ws.Cells.Item(8, 19).Value = 0.02099358841777628
$ws.Cells.Item(8, 20).Value = 0.02099358841777628

$ws.Cells.Item(9, 9).Value = 0.8065933093666526
$ws.Cells.Item(9, 10).Value = 0.8065933093666525
$ws.Cells.Item(9, 15).Value = 0.144012433133819
$ws.Cells.Item(9, 16).Value = 0.144012433133819
$ws.Cells.Item(9, 19).Value = 0.1161594650313509
$ws.Cells.Item(9, 20).Value = 0.1161594650313508

$ws.Cells.Item(10, 9).Value = 0.8065933093666526
$ws.Cells.Item(10, 10).Value = 0.8065933093666525
$ws.Cells.Item(10, 15).Value = 0.8299600903498424
$ws.Cells.Item(10, 16).Value = 0.8299600903498425
$ws.Cells.Item(10, 19).Value = 0.6694402559175254
$ws.Cells.Item(10, 20).Value = 0.6694402559175254
